# PRIXMM.xlsx — price list update + view-state refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price corrections on rows 45-46 (column D).
$ws.Range("D45").Value = 100
$ws.Range("D46").Value = 233

# Scroll the window so row 28 is the first visible row (topLeftCell = A28),
# then leave A2:B2 selected as the active selection, matching the sheet's
# last-saved view state.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2:B2").Select() | Out-Null
